$wb = $excel.ActiveWorkbook

# Sheet order (1-based) in this workbook:
# 1 Funciones_Objetivo
# 2 Restricciones_del_lider
# 3 Restricciones_del_follower
# 4 Punto_modificado
# 5 Vector_bf
# 6 Vector_BF   (name collides case-insensitively with "Vector_bf" -> address by index)
# 7 Vector_Alpha

# --- Restricciones_del_follower ---
$ws = $wb.Worksheets.Item(3)

# These columns hold numeric-looking values that must stay stored as TEXT
# (matching the regenerated source data), so mark them as Text before writing.
$ws.Range("B2:B4").NumberFormat = "@"
$ws.Range("D2:F4").NumberFormat = "@"

$ws.Range("A2").Value = "-9.75 + x + 1.2345679012345678y"
$ws.Range("B2").Value = "-4.25"
$ws.Range("D2").Value = "0.07"
$ws.Range("E2").Value = "0"
$ws.Range("F2").Value = "2.0"

$ws.Range("A3").Value = "-9.772 + x + 1.2400000000000002y"
$ws.Range("B3").Value = "7.772"
$ws.Range("D3").Value = "0.21"
$ws.Range("E3").Value = "6.800000000000001"
$ws.Range("F3").Value = "0"

$ws.Range("A4").Value = "-5.069999999999997 - 2x + 3.5802469135802464y"
$ws.Range("B4").Value = "3.9999999999999964"
$ws.Range("D4").Value = "0.4"
$ws.Range("E4").Value = "3.5"
$ws.Range("F4").Value = "5.8"

# --- Punto_modificado ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("A2:B2").NumberFormat = "@"
$ws.Range("A2").Value = "4.75"
$ws.Range("B2").Value = "4.05"

# --- Vector_bf ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "0.12108148148148201"

# --- Vector_BF (index 6) ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("A2:A3").NumberFormat = "@"
$ws.Range("A2").Value = "-3.3000000000000007"
$ws.Range("A3").Value = "-28.302864197530862"

# --- Vector_Alpha ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("A2").Value = 1.62
